$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.458.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.29%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.801.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.15%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "419.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.34%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.801.32"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.77%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.601"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.54%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.715"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.87%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.161"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -15.20%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000343"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -16.54%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.70%  "

# Row 14
$ws.Range("B14").Value = "Uniswap"
$ws.Range("C14").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +24.28%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.388.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.04%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "9.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.90%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.872.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.38%  "

# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.137"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.71%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.84%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.687.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.83%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.40%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "404.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.66%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.57%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.84%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.08%  "

# Row 26
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "36.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.18%  "

# Row 27
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.41%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.83%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.50%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "705.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.50%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.55%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.120"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.28%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.47%  "

# Row 35
$ws.Range("E35").Value = "  -9.56%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.00%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "37.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.77%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "54.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.80%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0757"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.41%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0450"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.39%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.97%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.25%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.60%  "

# Row 44
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.133"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.93%  "

# Row 45
$ws.Range("B45").Value = "LidoDAOToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.73%  "

# Row 46
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.26%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.10%  "

# Row 48
$ws.Range("E48").Value = "  -3.38%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.12%  "

# Row 50
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.49%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.62%  "
